$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Original layout:  A=bty_code  B=bty_iden  C=bty_labe  D=updated_at
# Target layout:    A=_airbyte_ab_id  B=_airbyte_emitted_at  C=bty_code
#                    D=bty_iden  E=bty_labe  F=_airbyte_additional_properties
#                    G=source_file_path  H=updated_at
# ---------------------------------------------------------------------------

# Insert two new columns at the front (for _airbyte_ab_id, _airbyte_emitted_at)
$ws.Range("A:B").Insert()

# After the first insert, the layout is:
#   C=bty_code  D=bty_iden  E=bty_labe  F=updated_at
# Insert two more columns right after E (bty_labe) for the Airbyte metadata
# columns that land between bty_labe and updated_at.
$ws.Range("F:G").Insert()

# Final layout:
#   A=_airbyte_ab_id  B=_airbyte_emitted_at  C=bty_code  D=bty_iden
#   E=bty_labe  F=_airbyte_additional_properties  G=source_file_path
#   H=updated_at

# --- Copy formatting from neighboring, already-correctly-styled cells so the
#     new columns reuse the existing cell styles instead of creating new ones.

# Header formatting (bold / centered / bordered) comes from the header row.
$ws.Range("C1").Copy() | Out-Null
$ws.Range("A1:B1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null

# Date-style formatting (the "YYYY-MM-DD HH:MM:SS" number format) comes from
# the updated_at column.
$ws.Range("H2").Copy() | Out-Null
$ws.Range("B2:B5").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false

# --- Header text ---
$ws.Range("A1").Value = "_airbyte_ab_id"
$ws.Range("B1").Value = "_airbyte_emitted_at"
$ws.Range("F1").Value = "_airbyte_additional_properties"
$ws.Range("G1").Value = "source_file_path"

# --- Data rows ---
$ids = @(
    "e15df2c2-1b47-4345-887c-40e02b52b201",
    "8ed9c208-27b0-4127-9e57-5b28ffff841c",
    "071ac219-4f4c-4bf4-805b-af36af18707d",
    "65f3601e-7c83-402d-bf86-f74ed903dbe4"
)
$sourcePath = "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/BIN_TYPE/2024_08_06_1722929004063_0.parquet"

for ($i = 0; $i -lt 4; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $ids[$i]
    $ws.Cells.Item($r, 2).Value = 45510.3079196875
    $ws.Cells.Item($r, 7).Value = $sourcePath
    $ws.Cells.Item($r, 8).Value = 45511.29381941652
}

$ws.Range("A1").Select() | Out-Null
